$d = $word.ActiveDocument

# Change 1: Insert "AL " before "{{NUMERO_VENDEDOR}}" in the seller notification phone clause.
$d.Content.Find.Execute("NOTIFICACIONES {{NUMERO_VENDEDOR}}", $false, $false, $false, $false, $false, $true, 1, $false, "NOTIFICACIONES AL {{NUMERO_VENDEDOR}}", 2)

# Change 2: Remove curly quotes around RELEVANTES.
$d.Content.Find.Execute("SEÑALADOS “RELEVANTES”, DONDE", $false, $false, $false, $false, $false, $true, 1, $false, "SEÑALADOS RELEVANTES, DONDE", 2)

# Changes 3 & 4: Replace the "PROMITENTE COMPRADOR" / "PROMITENTE VENDEDOR" signature
# line labels with sex-conditioned placeholders. Locate the specific paragraph that
# holds both labels (the tab-separated signature line near the end of the document)
# and scope the Find to just that paragraph's Range, so the earlier, unrelated
# "...PROMITENTE VENDEDOR LE ENTREGUE..." occurrence is left untouched, and so the
# tab-stop runs between the two labels are not disturbed by the replace.
$sigIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*PROMITENTE COMPRADOR*PROMITENTE VENDEDOR*") {
        $sigIndex = $i
        break
    }
}

if ($sigIndex -gt 0) {
    $sigRange1 = $d.Paragraphs.Item($sigIndex).Range
    $sigRange1.Find.Execute("COMPRADOR", $false, $false, $false, $false, $false, $true, 0, $false, "{{SEXO_4}}", 2)

    $sigRange2 = $d.Paragraphs.Item($sigIndex).Range
    $sigRange2.Find.Execute("VENDEDOR", $false, $false, $false, $false, $false, $true, 0, $false, "{{SEXO_2}}", 2)
} else {
    Write-Output "WARNING: signature paragraph not found"
}

Write-Output "done"
